$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old C/D header cells (columns C and D are no longer used)
$ws.Range("C1:D1").ClearContents()

# Column A: Classificatiecode (numeric 1..8), Column B: Onderdeel (names)
$ws.Range("A1").Value = "Classificatiecode"
$ws.Range("B1").Value = "Onderdeel"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Brandmeldinstallatie"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Sprinklerinstallatie"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Rolluik entree"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Rolluik fireshield"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Schuifdeur entree"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Inbraakbeveiliging"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Video bewakingssyteem"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Telefooninstallatie"

# New block in columns G/H: Eigenschappen / Waarde
$ws.Range("G1").Value = "Eigenschappen"
$ws.Range("H1").Value = "Waarde"

$ws.Range("G2").Value = "Aanwezig in Project"
$ws.Range("H2").Value = "True/false"

$ws.Range("G3").Value = "Standalone"
$ws.Range("H3").Value = "True/false"

# Column widths for new columns G and H (best-fit widths matching the source)
$ws.Columns.Item(7).ColumnWidth = 19.42578125
$ws.Columns.Item(8).ColumnWidth = 9.85546875
$ws.Columns.Item(7).BestFit = $true
$ws.Columns.Item(8).BestFit = $true

# Update the selected cell shown in the saved view
$ws.Range("K5").Select() | Out-Null
